# Applies RF007 - Gerenciar Avaliacoes text corrections (version 1.3 -> 1.4)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) "dos Avaliacoes cadastrados" -> "das Avaliacoes cadastradas"
$oldText1 = "SYSTEM exibe a listagem dos Avaliacoes cadastrados com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
$newText1 = "SYSTEM exibe a listagem das Avaliacoes cadastradas com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
foreach ($addr in @("D10", "D20", "D30", "D40", "D50", "D63", "D76", "D90")) {
    $cell = $ws.Range($addr)
    if ($cell.Text -eq $oldText1) {
        $cell.Value = $newText1
    }
}

# 2) "dos Avaliacoes com o Avaliacao excluido" -> "das Avaliacoes com a Avaliacao nao excluida"
$oldText2 = "SYSTEM exibe a listagem dos Avaliacoes com o Avaliacao excluido"
$newText2 = "SYSTEM exibe a listagem das Avaliacoes com a Avaliacao nao excluida"
$cell = $ws.Range("D13")
if ($cell.Text -eq $oldText2) {
    $cell.Value = $newText2
}

# 3) "dos Avaliacoes sem o Avaliacao excluido" -> "das Avaliacoes sem a Avaliacao excluida"
$oldText3 = "SYSTEM exibe a listagem dos Avaliacoes sem o Avaliacao excluido"
$newText3 = "SYSTEM exibe a listagem das Avaliacoes sem a Avaliacao excluida"
$cell = $ws.Range("D23")
if ($cell.Text -eq $oldText3) {
    $cell.Value = $newText3
}

# 4) "formulario para e alteracao de Avaliacao" -> "formulario para e alteracao da Avaliacao"
$oldText4 = "SYSTEM apresenta o formulario para e alteracao de Avaliacao"
$newText4 = "SYSTEM apresenta o formulario para e alteracao da Avaliacao"
foreach ($addr in @("D52", "D65", "D78")) {
    $cell = $ws.Range($addr)
    if ($cell.Text -eq $oldText4) {
        $cell.Value = $newText4
    }
}

# 5) "'Avaliado' somente leitura" -> "'Avaliado' estao em modo somente leitura"
$oldText5 = "Lider de Pessoas verifica que os campos 'Periodo Avaliativo', 'Perfil' e 'Avaliado' somente leitura"
$newText5 = "Lider de Pessoas verifica que os campos 'Periodo Avaliativo', 'Perfil' e 'Avaliado' estao em modo somente leitura"
foreach ($addr in @("B53", "B66", "B79")) {
    $cell = $ws.Range($addr)
    if ($cell.Text -eq $oldText5) {
        $cell.Value = $newText5
    }
}

# 6) "campos 'Nivel' preenchido corretamente" -> "campos 'Nivel' preenchidos corretamente"
$oldText6 = "SYSTEM apresenta em 'Metas' os campos 'Nivel' preenchido corretamente"
$newText6 = "SYSTEM apresenta em 'Metas' os campos 'Nivel' preenchidos corretamente"
$cell = $ws.Range("D95")
if ($cell.Text -eq $oldText6) {
    $cell.Value = $newText6
}
